# Add a new localization entry "PasswordResetUserNotFound" to the Türkçe
# localization workbook, as row 90 (A90/B90), matching the style of the
# existing data rows (A85:B85 is the last row where both cells share the
# "data row" style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from an existing fully-styled data row so the new row's
# cells pick up the same cellXfs (style) index as the rest of the table.
$ws.Range("A85:B85").Copy()
$ws.Range("A90:B90").PasteSpecial(-4122)

$ws.Range("A90").Value = "PasswordResetUserNotFound"
$ws.Range("B90").Value = "Bu e-posta adresine ait kullanici kaydi bulunmamistir. Lutfen e-posta adresinizin dogru oldugundan emin olunuz."
